# Generate Report for Handoff
#
# Updates the localization-status report:
#   * Refreshes the "Latest HO Xliff Generate Date" / "Latest Handoff
#     Datetime" timestamps for the 1bc11940-... row (handoff re-run).
#   * Fixes the "Priority" column for that same row (on both the zh-cn
#     and de-de sheets) from blank to "ht", matching the handoff type
#     that the mismatch message already calls out.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 10, 11, 12, 13)

# --- Overview sheet: "Latest HO Xliff Generate Date" (column G) ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-09-05 04:24:56"
}

# --- zh-cn sheet: "Latest Handoff Datetime" (column H) + "Priority" (column E) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("H$r").Value = "2016-09-05 04:24:50"
    $zhcn.Range("E$r").Value = "ht"
}

# --- de-de sheet: "Latest Handoff Datetime" (column H) + "Priority" (column E) ---
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("H$r").Value = "2016-09-05 04:24:56"
    $dede.Range("E$r").Value = "ht"
}
